# Auto-generated cell updates mirroring the cryptos-list refresh diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.089.90"
$ws.Range("E2").Value = "  +4.40%  "
$ws.Range("D3").Value = "1.906.38"
$ws.Range("E3").Value = "  +5.33%  "
$ws.Range("D4").Value = "'0.9999"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'251.84"
$ws.Range("E5").Value = "  +1.70%  "
$ws.Range("D6").Value = "'0.9996"
$ws.Range("E6").Value = "  -0.05%  "
$ws.Range("D7").Value = "'0.5094"
$ws.Range("E7").Value = "  +2.52%  "
$ws.Range("D8").Value = "'45.06"
$ws.Range("E8").Value = "  +4.33%  "
$ws.Range("D9").Value = "'0.3023"
$ws.Range("E9").Value = "  +8.73%  "
$ws.Range("D11").Value = "1.907.28"
$ws.Range("E11").Value = "  +5.36%  "
$ws.Range("D12").Value = "'17.26"
$ws.Range("D13").Value = "'0.07326"
$ws.Range("E13").Value = "  +3.57%  "
$ws.Range("D14").Value = "'0.6956"
$ws.Range("E14").Value = "  +7.73%  "
$ws.Range("D15").Value = "'86.58"
$ws.Range("E15").Value = "  +3.02%  "
$ws.Range("D16").Value = "'4.903"
$ws.Range("E16").Value = "  +4.51%  "
$ws.Range("D17").Value = "30.075.42"
$ws.Range("E17").Value = "  +4.36%  "
$ws.Range("D18").Value = "'0.000008162"
$ws.Range("E18").Value = "  +11.27%  "
$ws.Range("E19").Value = "  +0.05%  "
$ws.Range("D21").Value = "2.153.76"
$ws.Range("E21").Value = "  +5.37%  "
$ws.Range("D22").Value = "'0.9995"
$ws.Range("E22").Value = "  -0.01%  "
$ws.Range("D23").Value = "'4.821"
$ws.Range("E23").Value = "  +5.31%  "
$ws.Range("D24").Value = "'5.740"
$ws.Range("E24").Value = "  +7.54%  "
$ws.Range("D25").Value = "'9.267"
$ws.Range("E25").Value = "  +4.35%  "
$ws.Range("D26").Value = "'147.42"
$ws.Range("E26").Value = "  +3.57%  "
$ws.Range("D27").Value = "'134.86"
$ws.Range("E27").Value = "  +4.84%  "
$ws.Range("D28").Value = "'17.04"
$ws.Range("E28").Value = "  +4.21%  "
$ws.Range("D29").Value = "'1.995"
$ws.Range("E29").Value = "  +6.00%  "
$ws.Range("D30").Value = "'1.403"
$ws.Range("E30").Value = "  -0.75%  "
$ws.Range("D31").Value = "'4.247"
$ws.Range("E31").Value = "  +2.89%  "
$ws.Range("D32").Value = "'0.08812"
$ws.Range("E32").Value = "  +5.61%  "
$ws.Range("D33").Value = "'4.000"
$ws.Range("E33").Value = "  +5.15%  "
$ws.Range("D34").Value = "'0.05056"
$ws.Range("E34").Value = "  +2.16%  "
$ws.Range("D35").Value = "'1.137"
$ws.Range("E35").Value = "  +4.37%  "
$ws.Range("D36").Value = "'0.7186"
$ws.Range("E36").Value = "  +7.40%  "
$ws.Range("E37").Value = "  -0.69%  "
$ws.Range("E38").Value = "  +3.05%  "
$ws.Range("E39").Value = "  -0.24%  "
$ws.Range("D40").Value = "'0.9651"
$ws.Range("E40").Value = "  +1.11%  "
$ws.Range("E41").Value = "  +6.40%  "
$ws.Range("D42").Value = "'6.131"
$ws.Range("E42").Value = "  +0.70%  "
$ws.Range("D43").Value = "'0.4302"
$ws.Range("E43").Value = "  +5.48%  "
$ws.Range("D44").Value = "'104.68"
$ws.Range("E44").Value = "  +4.86%  "
$ws.Range("D45").Value = "'0.9990"
$ws.Range("E45").Value = "  -0.09%  "
$ws.Range("E46").Value = "  +6.49%  "
$ws.Range("D47").Value = "'0.1277"
$ws.Range("E47").Value = "  +4.89%  "
$ws.Range("E48").Value = "  +4.20%  "
$ws.Range("D49").Value = "'33.23"
$ws.Range("E49").Value = "  +5.24%  "
$ws.Range("D50").Value = "'8.404"
$ws.Range("E50").Value = "  +4.31%  "
$ws.Range("D51").Value = "'0.3810"
$ws.Range("E51").Value = "  +5.12%  "
